$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new employee row first (adds "Will" to shared strings before "Handle")
$ws.Range("A11").Value = 10
$ws.Range("A11").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B11").Value = "Will"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 1

# New header in E1
$ws.Range("E1").Value = "Handle"

# E2 gets its own (non-shared) formula
$ws.Range("E2").Formula = '="@"&B2'

# E3:E11 share a formula group
$ws.Range("E3:E11").Formula = '="@"&B3'

# Restore selection to mirror the author's final UI state as closely as possible
$null = $ws.Range("E2:E11").Select()
